$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 29 (Feria Lagunitas de Puerto Montt / Granada, Primera) and
# insert the copy above it, shifting rows 29-42 down to 30-43. This gives a
# new row 29 that starts out identical to the old row 29.
$ws.Rows.Item(29).Copy()
$ws.Rows.Item(29).Insert()

# Now adjust the new row 29 to reflect the new weekly data point.
$ws.Range("D29").Value = 44673
$ws.Range("M29").Value = 200
$ws.Range("Q29").Value = "$/caja 14 kilos empedrada"
$ws.Range("S29").Value = 1036
$ws.Range("T29").Value = 14
